# Add yesterday's and today's tasks to the progress report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: yesterday's task.
$ws.Cells.Item(3, 1).Value = 45752
$ws.Cells.Item(3, 2).Value = "Create github repo"
$ws.Cells.Item(3, 3).Value = "Complete"
$ws.Cells.Item(3, 4).Value = "Adam Rodi"

# Row 2 existing task status changes from "In Progress" to "Complete".
$ws.Range("C2").Value = "Complete"

# Row 4: today's task.
$ws.Cells.Item(4, 1).Value = 45753
$ws.Cells.Item(4, 2).Value = "Decide on project idea"
$ws.Cells.Item(4, 3).Value = "In Progress"
$ws.Cells.Item(4, 4).Value = "Everyone"

# Row 5: today's task.
$ws.Cells.Item(5, 1).Value = 45753
$ws.Cells.Item(5, 2).Value = "Design module communication diagram"
$ws.Cells.Item(5, 3).Value = "In Progress"
$ws.Cells.Item(5, 4).Value = "Everyone"

# Copy style (date format, center alignment) from A2 down to the new date cells.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Copy style from B2/C2/D2 down to the new rows so formatting matches.
$ws.Range("B2:D2").Copy() | Out-Null
$ws.Range("B3:D5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Adjust column widths to fit the new content (closest achievable values
# given the host's internal pixel-snapped width granularity).
$ws.Columns.Item(2).ColumnWidth = 31.75
$ws.Columns.Item(4).ColumnWidth = 10.75

# Update selection to match the recorded state after editing.
$ws.Range("B6").Select() | Out-Null
